$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 570
$ws.Range("F2").Value = 368
$ws.Range("G2").Value = 369

$ws.Range("D3").Value = 585
$ws.Range("F3").Value = 370
$ws.Range("G3").Value = 371

$ws.Range("D4").Value = 586
$ws.Range("F4").Value = 235
$ws.Range("G4").Value = 236

$ws.Range("D5").Value = 588
$ws.Range("F5").Value = 281
$ws.Range("G5").Value = 284

$ws.Range("D6").Value = 587
$ws.Range("F6").Value = 490
$ws.Range("G6").Value = 489

$ws.Range("D7").Value = 589
$ws.Range("F7").Value = 283
$ws.Range("G7").Value = 286

$ws.Range("D8").Value = 590
$ws.Range("F8").Value = 233
$ws.Range("G8").Value = 234

$ws.Range("D9").Value = 591
$ws.Range("F9").Value = 372
$ws.Range("G9").Value = 373

$ws.Range("D10").Value = 592
$ws.Range("F10").Value = 231
$ws.Range("G10").Value = 232

$ws.Range("D11").Value = 593
$ws.Range("F11").Value = 383
$ws.Range("G11").Value = 384

$ws.Range("D12").Value = 594
$ws.Range("F12").Value = 391
$ws.Range("G12").Value = 392

$ws.Range("D13").Value = 596
$ws.Range("F13").Value = 433
$ws.Range("G13").Value = 436

$ws.Range("D14").Value = 595
$ws.Range("F14").Value = 434
$ws.Range("G14").Value = 435

$ws.Range("D15").Value = 597
$ws.Range("F15").Value = 410
$ws.Range("G15").Value = 408

$ws.Range("D16").Value = 598
$ws.Range("F16").Value = 411
$ws.Range("G16").Value = 409

$ws.Range("D17").Value = 599
$ws.Range("F17").Value = 348
$ws.Range("G17").Value = 347

$ws.Range("C18").Value = 620
$ws.Range("D18").Value = 621
$ws.Range("F18").Value = 622
$ws.Range("G18").Value = 623

$ws.Range("C19").Value = 552
$ws.Range("D19").Value = 601
$ws.Range("F19").Value = 559
$ws.Range("G19").Value = 566

$ws.Range("C20").Value = 551
$ws.Range("D20").Value = 602
$ws.Range("F20").Value = 558
$ws.Range("G20").Value = 565

$ws.Range("C21").Value = 631
$ws.Range("D21").Value = 632
$ws.Range("F21").Value = 639
$ws.Range("G21").Value = 640

$ws.Range("C22").Value = 550
$ws.Range("D22").Value = 603
$ws.Range("F22").Value = 557
$ws.Range("G22").Value = 564

$ws.Range("C23").Value = 630
$ws.Range("D23").Value = 633
$ws.Range("F23").Value = 638
$ws.Range("G23").Value = 641

$ws.Range("C24").Value = 549
$ws.Range("D24").Value = 604
$ws.Range("F24").Value = 556
$ws.Range("G24").Value = 563

$ws.Range("D25").Value = 613
$ws.Range("F25").Value = 277
$ws.Range("G25").Value = 278

$ws.Range("D26").Value = 612
$ws.Range("F26").Value = 279
$ws.Range("G26").Value = 280

$ws.Range("G27").Select() | Out-Null
